$wb = $excel.ActiveWorkbook

# Sheet "qtd": rows 21-23 (MA..RO tied block) re-ranked TO,MS,MT -> MS,MT,TO
$wsQtd = $wb.Worksheets.Item("qtd")
$wsQtd.Range("A21").Value = "MS"
$wsQtd.Range("A22").Value = "MT"
$wsQtd.Range("A23").Value = "TO"

# Sheet "max-arrecad": rows 25-26 (AC,TO tied block) swapped
$wsMax = $wb.Worksheets.Item("max-arrecad")
$wsMax.Range("A25").Value = "TO"
$wsMax.Range("A26").Value = "AC"

# Sheet "tx-sucesso": rows 23-25 (AC,RO,SE tied block at 0) re-ranked AC,RO,SE -> RO,SE,AC
$wsTx = $wb.Worksheets.Item("tx-sucesso")
$wsTx.Range("A23").Value = "RO"
$wsTx.Range("A24").Value = "SE"
$wsTx.Range("A25").Value = "AC"
